# Updated cryptos list on Wed Jan 31 06:35:10 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for the existing coins,
# and rotates out "ApeXProtocol" (old row 43) so every coin from Maker (old
# row 44) down to RocketPoolETH (old row 51) shifts up one row, with
# THORChain appended as the new row 51.
#
# Note: several Price values look like plain numbers (e.g. "306.70",
# "0.0800"). Assigning those bare would let Excel auto-convert them to
# numeric values and silently drop significant trailing zeros, which would
# not match the source text cells. Prefixing with a leading single-quote
# forces Excel to keep them as literal text, matching the original
# (inline-string) cell contents exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.030.83'
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").Value = '2.338.87'
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''306.70'
$ws.Range("E5").Value = '  -1.44%  '

$ws.Range("D6").Value = '''101.09'
$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("E7").Value = '  -3.43%  '

$ws.Range("D9").Value = '''0.511'
$ws.Range("E9").Value = '  -3.37%  '

$ws.Range("D10").Value = '''34.97'
$ws.Range("E10").Value = '  -1.92%  '

$ws.Range("D11").Value = '''52.30'
$ws.Range("E11").Value = '  +1.12%  '

$ws.Range("D12").Value = '''0.0800'
$ws.Range("E12").Value = '  -1.34%  '

$ws.Range("E13").Value = '  +0.90%  '

$ws.Range("E14").Value = '  -1.73%  '

$ws.Range("D15").Value = '''15.93'
$ws.Range("E15").Value = '  +6.27%  '

$ws.Range("D16").Value = '2.296.36'
$ws.Range("E16").Value = '  +0.62%  '

$ws.Range("D17").Value = '''0.814'
$ws.Range("E17").Value = '  +0.92%  '

$ws.Range("D18").Value = '42.960.30'
$ws.Range("E18").Value = '  -0.70%  '

$ws.Range("D19").Value = '''11.81'
$ws.Range("E19").Value = '  -3.90%  '

$ws.Range("E21").Value = '  -2.20%  '

$ws.Range("D22").Value = '''67.91'
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("D23").Value = '''237.05'
$ws.Range("E23").Value = '  -1.75%  '

$ws.Range("D24").Value = '''2.03'
$ws.Range("E24").Value = '  +1.16%  '

$ws.Range("E25").Value = '  -2.11%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").Value = '''25.49'
$ws.Range("E27").Value = '  +3.02%  '

$ws.Range("E29").Value = '  +5.76%  '

$ws.Range("D30").Value = '''35.07'
$ws.Range("E30").Value = '  -4.44%  '

$ws.Range("E31").Value = '  -2.36%  '

$ws.Range("D32").Value = '''163.56'
$ws.Range("E32").Value = '  -3.11%  '

$ws.Range("E33").Value = '  -0.09%  '

$ws.Range("E34").Value = '  -2.68%  '

$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("E36").Value = '  -2.51%  '

$ws.Range("D37").Value = '''4.62'
$ws.Range("E37").Value = '  +6.16%  '

$ws.Range("D38").Value = '''0.0728'
$ws.Range("E38").Value = '  -1.83%  '

$ws.Range("E39").Value = '  -1.36%  '

$ws.Range("D40").Value = '''2.93'
$ws.Range("E40").Value = '  -4.38%  '

$ws.Range("E41").Value = '  -3.60%  '

$ws.Range("E42").Value = '  -2.00%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.007.28'
$ws.Range("E43").Value = '  +2.12%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0286'
$ws.Range("E44").Value = '  -0.99%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''18.69'
$ws.Range("E45").Value = '  -3.77%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '''10.20'
$ws.Range("E46").Value = '  +3.24%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''2.94'
$ws.Range("E47").Value = '  -1.67%  '

$ws.Range("B48").Value = 'MultiversX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D48").Value = '''56.10'
$ws.Range("E48").Value = '  +1.28%  '

$ws.Range("B49").Value = 'HuobiToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D49").Value = '''2.89'
$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.564.08'
$ws.Range("E50").Value = '  +1.20%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = '''4.71'
$ws.Range("E51").Value = '  +2.74%  '
